$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routes")

$ws.Range("A20").Value = "EuropeanWings"
$ws.Range("C20").Value = "LFOB"
$ws.Range("B20").Value = "Paris-Beauvais-Tille"
$ws.Range("D20").Value = "Hungary-Budapest-Listz"
$ws.Range("E20").Value = "LHBP"

$ws.Range("D20").VerticalAlignment = -4108

$ws.Range("E20").Select()
